# Updated cryptos list — applies refreshed price/volume figures to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D occasionally hold numeric-looking text (e.g. "215.43").
# Force those specific cells to Text format first so Excel keeps them as
# strings instead of silently converting them to numbers.
$textCells = @('D5', 'D9', 'D17', 'D20', 'D22', 'D25', 'D29', 'D30', 'D31', 'D34', 'D36', 'D42', 'D47', 'D48')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = '@'
}

$ws.Range('D2').Value = '25.947.41'
$ws.Range('E2').Value = '  -1.28%  '
$ws.Range('D3').Value = '1.637.80'
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('D5').Value = '215.43'
$ws.Range('E5').Value = '  -0.88%  '
$ws.Range('E6').Value = '  -0.44%  '
$ws.Range('E7').Value = '  +0.30%  '
$ws.Range('E8').Value = '  -0.72%  '
$ws.Range('D9').Value = '0.0638'
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('E10').Value = '  -1.95%  '
$ws.Range('D12').Value = '1.865.23'
$ws.Range('E12').Value = '  -0.52%  '
$ws.Range('E13').Value = '  -0.99%  '
$ws.Range('D14').Value = '1.636.53'
$ws.Range('E14').Value = '  -0.67%  '
$ws.Range('E15').Value = '  -1.28%  '
$ws.Range('E16').Value = '  -0.21%  '
$ws.Range('D17').Value = '62.99'
$ws.Range('E17').Value = '  -0.88%  '
$ws.Range('D18').Value = '25.990.95'
$ws.Range('E19').Value = '  +0.28%  '
$ws.Range('D20').Value = '193.07'
$ws.Range('E20').Value = '  -1.78%  '
$ws.Range('E21').Value = '  -1.74%  '
$ws.Range('D22').Value = '9.94'
$ws.Range('E22').Value = '  -1.41%  '
$ws.Range('E23').Value = '  -0.92%  '
$ws.Range('E24').Value = '  +0.76%  '
$ws.Range('D25').Value = '143.99'
$ws.Range('E25').Value = '  +0.52%  '
$ws.Range('E26').Value = '  +0.28%  '
$ws.Range('E27').Value = '  +3.19%  '
$ws.Range('E28').Value = '  -1.91%  '
$ws.Range('D29').Value = '15.58'
$ws.Range('D30').Value = '1.24'
$ws.Range('E30').Value = '  -0.83%  '
$ws.Range('D31').Value = '0.0503'
$ws.Range('E31').Value = '  -0.72%  '
$ws.Range('E32').Value = '  -1.24%  '
$ws.Range('E33').Value = '  -0.46%  '
$ws.Range('D34').Value = '1.54'
$ws.Range('E34').Value = '  -4.67%  '
$ws.Range('E35').Value = '  +1.63%  '
$ws.Range('D36').Value = '0.900'
$ws.Range('E36').Value = '  -1.50%  '
$ws.Range('D37').Value = '1.136.25'
$ws.Range('E37').Value = '  -0.37%  '
$ws.Range('E38').Value = '  -1.95%  '
$ws.Range('E39').Value = '  -1.45%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('D42').Value = '5.49'
$ws.Range('E42').Value = '  -3.95%  '
$ws.Range('E43').Value = '  -0.99%  '
$ws.Range('E44').Value = '  -0.49%  '
$ws.Range('D45').Value = '1.774.89'
$ws.Range('E46').Value = '  +2.12%  '
$ws.Range('D47').Value = '56.76'
$ws.Range('E47').Value = '  +0.56%  '
$ws.Range('D48').Value = '0.0530'
$ws.Range('E48').Value = '  +3.19%  '
$ws.Range('E49').Value = '  +0.59%  '
$ws.Range('E50').Value = '  -0.36%  '
$ws.Range('E51').Value = '  -0.68%  '
